$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("C2").Value = "83.1102981641134"
$ws.Range("I2").Value = "-476.85244136181313"
$ws.Range("C3").Value = "400.0000000176049"
$ws.Range("I3").Value = "-462.8409244919828"
$ws.Range("C4").Value = "400.0000000356312"
$ws.Range("I4").Value = "-390.2720508695672"
$ws.Range("C5").Value = "420.0000000374189"
$ws.Range("I5").Value = "-486.8191586029397"
$ws.Range("I6").Value = "-138.22243104805233"
$ws.Range("C7").Value = "358.07549512951493"
$ws.Range("I7").Value = "37.49204221702779"
$ws.Range("C8").Value = "299.3835951197538"
$ws.Range("I8").Value = "-119.60442139538192"
$ws.Range("C9").Value = "400.0000000410847"
$ws.Range("C10").Value = "400.00000004108404"
$ws.Range("I10").Value = "675.7517231211192"
$ws.Range("J10").Value = "796.7433714822"
$ws.Range("C11").Value = "400.00000004041436"
$ws.Range("I11").Value = "1264.0717322700705"
$ws.Range("C12").Value = "400.0000000412478"
$ws.Range("I12").Value = "1124.264125978972"
$ws.Range("C13").Value = "400.0000000410847"
$ws.Range("I13").Value = "800.37976258646"
$ws.Range("C14").Value = "400.0000000412478"
$ws.Range("I14").Value = "-18.858220925469503"
$ws.Range("C15").Value = "489.62065399352576"
$ws.Range("I15").Value = "959.0322788555521"
$ws.Range("C16").Value = "400.0000001247187"
$ws.Range("C17").Value = "400.0000001247187"
$ws.Range("C18").Value = "400.0000001250004"
$ws.Range("I18").Value = "253.87520101539747"
$ws.Range("C19").Value = "400.00000012471844"
$ws.Range("I19").Value = "-581.1465953978322"
$ws.Range("C20").Value = "500.0"
$ws.Range("I20").Value = "-0.019818332195148107"
$ws.Range("C21").Value = "170.39999387651875"
$ws.Range("I21").Value = "-299.487372693519"
$ws.Range("C22").Value = "400.00000011271607"
$ws.Range("I22").Value = "-732.4069906486361"
$ws.Range("C23").Value = "400.000000119644"
$ws.Range("I23").Value = "-806.047954553181"
$ws.Range("C24").Value = "399.99999994878533"
$ws.Range("I24").Value = "-501.27656027536034"
$ws.Range("C25").Value = "420.00000009118696"
$ws.Range("I25").Value = "-138.22243104805236"
$ws.Range("C26").Value = "489.6201370111077"
$ws.Range("C27").Value = "400.000000119644"
$ws.Range("I27").Value = "-119.60442139538192"
$ws.Range("I29").Value = "58.94824951638855"
$ws.Range("I31").Value = "-62.48896717923746"
$ws.Range("I32").Value = "-62.48896717923779"
$ws.Range("I33").Value = "203.2570807097871"
$ws.Range("I34").Value = "203.25708070978746"
$ws.Range("I35").Value = "44.3725497239757"
$ws.Range("I36").Value = "44.37254972397571"

$ws = $wb.Worksheets.Item(2)
$ws.Range("C2").Value = "83.1102981641134"
$ws.Range("I2").Value = "-438.27835994496365"
$ws.Range("I3").Value = "-462.80828808309724"
$ws.Range("C4").Value = "400.0000000356312"
$ws.Range("I4").Value = "-389.44181835079087"
$ws.Range("C5").Value = "420.0000000374189"
$ws.Range("I5").Value = "-486.8009536421274"
$ws.Range("C6").Value = "456.7606477280851"
$ws.Range("I6").Value = "-119.07707149006714"
$ws.Range("C7").Value = "298.9999999997249"
$ws.Range("I7").Value = "37.557303398225415"
$ws.Range("C8").Value = "299.3835951197538"
$ws.Range("I8").Value = "-119.1950950003588"
$ws.Range("C9").Value = "400.0000000410847"
$ws.Range("C10").Value = "400.00000004108404"
$ws.Range("I10").Value = "694.268708738456"
$ws.Range("J10").Value = "796.7433714822"
$ws.Range("C11").Value = "400.00000004041436"
$ws.Range("I11").Value = "1264.336978268255"
$ws.Range("C12").Value = "400.0000000412478"
$ws.Range("I12").Value = "1125.5689278971631"
$ws.Range("C13").Value = "400.0000000410847"
$ws.Range("I13").Value = "800.3900951422893"
$ws.Range("C14").Value = "400.0000000412478"
$ws.Range("C15").Value = "491.7523107815771"
$ws.Range("I15").Value = "959.2364766498399"
$ws.Range("C16").Value = "400.0000001247187"
$ws.Range("C17").Value = "400.0000001247187"
$ws.Range("I17").Value = "-35.768021113128384"
$ws.Range("C18").Value = "400.0000001250004"
$ws.Range("I18").Value = "253.94594162461001"
$ws.Range("C19").Value = "400.00000012471844"
$ws.Range("I19").Value = "-579.7841364209398"
$ws.Range("C20").Value = "497.78327147087157"
$ws.Range("I20").Value = "-0.025255159646968898"
$ws.Range("C21").Value = "170.39999387651875"
$ws.Range("I21").Value = "-335.3871663101083"
$ws.Range("C22").Value = "400.00000011271607"
$ws.Range("I22").Value = "-732.4391777331941"
$ws.Range("C23").Value = "400.000000119644"
$ws.Range("I23").Value = "-806.8760812053656"
$ws.Range("C24").Value = "399.99999994878533"
$ws.Range("I24").Value = "-501.2947652361671"
$ws.Range("C25").Value = "420.00000009118696"
$ws.Range("I25").Value = "-119.07707149006718"
$ws.Range("C26").Value = "491.75177291206927"
$ws.Range("C27").Value = "400.000000119644"
$ws.Range("I27").Value = "-119.1950950003588"
$ws.Range("I29").Value = "31.639467094110003"
$ws.Range("I30").Value = "31.639467094109992"
$ws.Range("I31").Value = "-63.828026465895995"
$ws.Range("I32").Value = "-63.82802646589035"
$ws.Range("I33").Value = "204.55163692264557"
$ws.Range("I34").Value = "204.55163692264594"
$ws.Range("I35").Value = "45.43150729770366"
$ws.Range("I36").Value = "45.431507297703654"

$ws = $wb.Worksheets.Item(3)
$ws.Range("C2").Value = "83.1102981641134"
$ws.Range("I2").Value = "-43.459327353400724"
$ws.Range("C3").Value = "400.0000000176049"
$ws.Range("I3").Value = "-462.8898979642314"
$ws.Range("C4").Value = "400.0000000356312"
$ws.Range("I4").Value = "-390.88486197237575"
$ws.Range("C5").Value = "420.0000000374189"
$ws.Range("I5").Value = "-486.9553120373505"
$ws.Range("C6").Value = "452.22080960748644"
$ws.Range("I6").Value = "77.98291472269224"
$ws.Range("C7").Value = "298.9999999997249"
$ws.Range("I7").Value = "-3.752525846700998e-05"
$ws.Range("I8").Value = "-119.95826046392548"
$ws.Range("C9").Value = "400.0000000410847"
$ws.Range("I9").Value = "9.631388136617925"
$ws.Range("C10").Value = "400.00000004108404"
$ws.Range("I10").Value = "695.7389085759012"
$ws.Range("J10").Value = "796.7433714822"
$ws.Range("C11").Value = "400.00000004041436"
$ws.Range("I11").Value = "1263.6777624668707"
$ws.Range("C12").Value = "400.0000000412478"
$ws.Range("I12").Value = "1122.5789147601035"
$ws.Range("C13").Value = "400.0000000410847"
$ws.Range("I13").Value = "799.8954017744701"
$ws.Range("C14").Value = "400.0000000412478"
$ws.Range("C15").Value = "487.2166134412663"
$ws.Range("I15").Value = "958.7280351725541"
$ws.Range("C16").Value = "400.0000001247187"
$ws.Range("I16").Value = "-0.1583245037458638"
$ws.Range("C17").Value = "400.0000001247187"
$ws.Range("I17").Value = "-348.78939788391824"
$ws.Range("C18").Value = "400.0000001250004"
$ws.Range("I18").Value = "253.7652409305541"
$ws.Range("C19").Value = "400.00000012471844"
$ws.Range("I19").Value = "-581.8224010630845"
$ws.Range("C20").Value = "411.58897700967134"
$ws.Range("C21").Value = "170.39999387651875"
$ws.Range("I21").Value = "-649.6806607976703"
$ws.Range("C22").Value = "400.00000011271607"
$ws.Range("I22").Value = "-732.3586913433825"
$ws.Range("C23").Value = "400.000000119644"
$ws.Range("I23").Value = "-805.4366978323034"
$ws.Range("C24").Value = "399.99999994878533"
$ws.Range("I24").Value = "-501.1404068406879"
$ws.Range("C25").Value = "420.00000009118696"
$ws.Range("I25").Value = "77.98291472269227"
$ws.Range("C26").Value = "487.2157974776127"
$ws.Range("I26").Value = "37.39270553515247"
$ws.Range("C27").Value = "400.000000119644"
$ws.Range("I27").Value = "-119.95826046392548"
$ws.Range("I29").Value = "46.592598269624354"
$ws.Range("I31").Value = "-103.59666999025876"
$ws.Range("I32").Value = "-103.59666999025876"
$ws.Range("I33").Value = "248.12815062631546"
$ws.Range("I34").Value = "248.12815062631594"
$ws.Range("I35").Value = "38.59861173073795"
$ws.Range("I36").Value = "38.598611730738355"

$ws = $wb.Worksheets.Item(4)
$ws.Range("C2").Value = "83.1102981641134"
$ws.Range("I2").Value = "-200.02805901586385"
$ws.Range("C3").Value = "400.0000000176049"
$ws.Range("I3").Value = "-188.72741634293783"
$ws.Range("C4").Value = "400.0000000356312"
$ws.Range("I4").Value = "-390.3225852764998"
$ws.Range("C5").Value = "420.0000000374189"
$ws.Range("I5").Value = "-488.61484594877106"
$ws.Range("C6").Value = "396.7433828793597"
$ws.Range("C7").Value = "308.9473409650441"
$ws.Range("I7").Value = "148.18953600778673"
$ws.Range("C8").Value = "299.3835951197538"
$ws.Range("I8").Value = "-119.62521800637376"
$ws.Range("C9").Value = "400.0000000410847"
$ws.Range("I9").Value = "79.5455848850956"
$ws.Range("C10").Value = "400.00000004108404"
$ws.Range("I10").Value = "796.3941013969464"
$ws.Range("J10").Value = "796.7433714822"
$ws.Range("C11").Value = "400.00000004041436"
$ws.Range("I11").Value = "1332.2403989625527"
$ws.Range("C12").Value = "400.0000000412478"
$ws.Range("I12").Value = "1124.2458387738018"
$ws.Range("C13").Value = "400.0000000410847"
$ws.Range("I13").Value = "795.7007037537583"
$ws.Range("C14").Value = "400.0000000412478"
$ws.Range("I14").Value = "96.39559408059051"
$ws.Range("C15").Value = "489.5607066648654"
$ws.Range("I15").Value = "1034.3929433734363"
$ws.Range("C16").Value = "400.0000001247187"
$ws.Range("I16").Value = "-1.8282474118247276"
$ws.Range("I17").Value = "-91.22133875848775"
$ws.Range("C18").Value = "400.0000001250004"
$ws.Range("I18").Value = "347.8940392377415"
$ws.Range("C19").Value = "400.00000012471844"
$ws.Range("I19").Value = "-581.2579688043515"
$ws.Range("C20").Value = "500.0"
$ws.Range("C21").Value = "170.39999387651875"
$ws.Range("I21").Value = "-391.8033047288304"
$ws.Range("C22").Value = "400.00000011271607"
$ws.Range("I22").Value = "-607.9835232735547"
$ws.Range("C23").Value = "400.000000119644"
$ws.Range("I23").Value = "-805.997548325663"
$ws.Range("C24").Value = "399.99999994878533"
$ws.Range("I24").Value = "-499.48087292952425"
$ws.Range("C25").Value = "420.00000009118696"
$ws.Range("I25").Value = "1.139651416505324e-05"
$ws.Range("C26").Value = "489.5687839147222"
$ws.Range("I26").Value = "148.18953600778673"
$ws.Range("C27").Value = "400.000000119644"
$ws.Range("I27").Value = "-119.62521800637376"
$ws.Range("I28").Value = "79.54558488509558"
$ws.Range("I29").Value = "-95.97132749109326"
$ws.Range("I31").Value = "-154.32492248856872"
$ws.Range("I32").Value = "-154.32492248856872"
$ws.Range("I33").Value = "371.560017721006"
$ws.Range("I34").Value = "371.5600177210063"
$ws.Range("I35").Value = "-32.965799263576045"
$ws.Range("I36").Value = "-32.965799263576045"

$ws = $wb.Worksheets.Item(5)
$ws.Range("C2").Value = "83.1102981641134"
$ws.Range("I2").Value = "-44.8992741717568"
$ws.Range("C3").Value = "400.0000000176049"
$ws.Range("I3").Value = "-468.3658101196141"
$ws.Range("C4").Value = "400.0000000356312"
$ws.Range("I4").Value = "-389.9585527545011"
$ws.Range("C5").Value = "420.0000000374189"
$ws.Range("I5").Value = "-488.8470349463189"
$ws.Range("C6").Value = "455.44341753259465"
$ws.Range("I6").Value = "77.5734253890989"
$ws.Range("C7").Value = "298.9999999997249"
$ws.Range("I7").Value = "26.115655347194654"
$ws.Range("C8").Value = "299.3835951197538"
$ws.Range("I8").Value = "-119.44871549219265"
$ws.Range("C9").Value = "400.0000000410847"
$ws.Range("I9").Value = "84.20643664574385"
$ws.Range("C10").Value = "400.00000004108404"
$ws.Range("I10").Value = "699.7787403043855"
$ws.Range("J10").Value = "796.7433714822"
$ws.Range("C11").Value = "400.00000004041436"
$ws.Range("I11").Value = "1221.426143416667"
$ws.Range("C12").Value = "400.0000000412478"
$ws.Range("I12").Value = "1124.7715720294748"
$ws.Range("C13").Value = "400.0000000410847"
$ws.Range("I13").Value = "795.0711089750715"
$ws.Range("C14").Value = "400.0000000412478"
$ws.Range("I14").Value = "9.351550134491621e-05"
$ws.Range("C15").Value = "490.442712809505"
$ws.Range("I15").Value = "925.9076053495214"
$ws.Range("C16").Value = "400.0000001247187"
$ws.Range("I16").Value = "-2.0784941312163503"
$ws.Range("C17").Value = "400.0000001247187"
$ws.Range("I17").Value = "-204.36078886408015"
$ws.Range("I18").Value = "240.61113024225654"
$ws.Range("C19").Value = "400.00000012471844"
$ws.Range("I19").Value = "-580.6386903821603"
$ws.Range("C20").Value = "471.16514429229125"
$ws.Range("C21").Value = "170.39999387651875"
$ws.Range("I21").Value = "-505.5129240445194"
$ws.Range("C22").Value = "400.00000011271607"
$ws.Range("I22").Value = "-726.9581639266507"
$ws.Range("C23").Value = "400.000000119644"
$ws.Range("I23").Value = "-806.360657487145"
$ws.Range("C24").Value = "399.99999994878533"
$ws.Range("I24").Value = "-499.2486839319266"
$ws.Range("C25").Value = "420.00000009118696"
$ws.Range("I25").Value = "77.5734253890976"
$ws.Range("C26").Value = "490.45148875204603"
$ws.Range("C27").Value = "400.000000119644"
$ws.Range("I27").Value = "-119.44871549219265"
$ws.Range("I28").Value = "84.20643664574385"
$ws.Range("I29").Value = "-51.38099751579956"
$ws.Range("I30").Value = "-51.38099751579956"
$ws.Range("I31").Value = "-174.6723369079266"
$ws.Range("I32").Value = "-174.67233690792656"
$ws.Range("I33").Value = "395.7167707360286"
$ws.Range("I34").Value = "395.7167707360287"
$ws.Range("I35").Value = "-36.883906624453886"
$ws.Range("I36").Value = "-36.88390662445388"

$ws = $wb.Worksheets.Item(6)
$ws.Range("C2").Value = "83.1102981641134"
$ws.Range("I2").Value = "-43.455798537492015"
$ws.Range("C3").Value = "400.0000000176049"
$ws.Range("I3").Value = "-462.7978105465429"
$ws.Range("C4").Value = "400.0000000356312"
$ws.Range("I4").Value = "-541.1365909704335"
$ws.Range("C5").Value = "420.0000000374189"
$ws.Range("I5").Value = "-436.82553705328877"
$ws.Range("C6").Value = "457.1425904034508"
$ws.Range("I6").Value = "78.38586389296611"
$ws.Range("C7").Value = "307.1591138469874"
$ws.Range("I7").Value = "37.57857936718954"
$ws.Range("C8").Value = "299.3835951197538"
$ws.Range("I8").Value = "-194.7014873295867"
$ws.Range("C9").Value = "400.0000000410847"
$ws.Range("C10").Value = "400.00000004108404"
$ws.Range("I10").Value = "699.8503895156998"
$ws.Range("J10").Value = "796.7433714822"
$ws.Range("I11").Value = "1264.4204343226447"
$ws.Range("C12").Value = "400.0000000412478"
$ws.Range("I12").Value = "878.0615007629301"
$ws.Range("C13").Value = "400.0000000410847"
$ws.Range("I13").Value = "750.622473841529"
$ws.Range("C14").Value = "400.0000000412478"
$ws.Range("C15").Value = "492.1535966624448"
$ws.Range("I15").Value = "959.3009700570339"
$ws.Range("C16").Value = "400.0000001247187"
$ws.Range("I16").Value = "-49.67016039426136"
$ws.Range("C17").Value = "400.0000001247187"
$ws.Range("I17").Value = "-333.02251388099353"
$ws.Range("C18").Value = "400.0000001250004"
$ws.Range("I18").Value = "253.9697685109613"
$ws.Range("C19").Value = "400.00000012471844"
$ws.Range("I19").Value = "-824.6962847956155"
$ws.Range("I20").Value = "-49.760060786513996"
$ws.Range("C21").Value = "170.39999387651875"
$ws.Range("I21").Value = "-635.363972920083"
$ws.Range("C22").Value = "400.00000011271607"
$ws.Range("I22").Value = "-732.4495110296849"
$ws.Range("C23").Value = "400.000000119644"
$ws.Range("I23").Value = "-655.5660790526186"
$ws.Range("C24").Value = "399.99999994878533"
$ws.Range("I24").Value = "-551.2701818250248"
$ws.Range("C25").Value = "420.00000009118696"
$ws.Range("I25").Value = "78.38586389296617"
$ws.Range("C26").Value = "492.1431057592792"
$ws.Range("C27").Value = "400.000000119644"
$ws.Range("I27").Value = "-194.70148732959808"
$ws.Range("I28").Value = "-45.54430471680575"
$ws.Range("I29").Value = "-119.49404641389964"
$ws.Range("I30").Value = "-119.49404641389962"
$ws.Range("I31").Value = "-110.75109843305171"
$ws.Range("I32").Value = "-110.75109843305171"
$ws.Range("I33").Value = "203.4786096279871"
$ws.Range("I34").Value = "203.47860962798728"
$ws.Range("I35").Value = "93.12605090606183"
$ws.Range("I36").Value = "93.1260509060618"

$ws = $wb.Worksheets.Item(7)
$ws.Range("C2").Value = "83.1102981641134"
$ws.Range("I2").Value = "-200.17285119375558"
$ws.Range("C3").Value = "400.0000000176049"
$ws.Range("I3").Value = "-462.9180432847205"
$ws.Range("C4").Value = "400.0000000356312"
$ws.Range("I4").Value = "-389.53312898750676"
$ws.Range("C5").Value = "420.0000000374189"
$ws.Range("I5").Value = "-217.1338208954108"
$ws.Range("C6").Value = "396.7433714646718"
$ws.Range("C7").Value = "361.0749511408736"
$ws.Range("C8").Value = "299.3835951197538"
$ws.Range("I8").Value = "-119.22099286304184"
$ws.Range("C9").Value = "400.0000000410847"
$ws.Range("C10").Value = "400.00000004108404"
$ws.Range("I10").Value = "796.3941013969464"
$ws.Range("J10").Value = "796.7433714822"
$ws.Range("C11").Value = "400.00000004041436"
$ws.Range("I11").Value = "1263.4892597416108"
$ws.Range("C12").Value = "400.0000000412478"
$ws.Range("I12").Value = "1125.6969459480988"
$ws.Range("C13").Value = "400.0000000410847"
$ws.Range("I13").Value = "924.8714505300187"
$ws.Range("C14").Value = "400.0000000412478"
$ws.Range("I14").Value = "96.14590786913759"
$ws.Range("C15").Value = "491.82954861820764"
$ws.Range("I15").Value = "958.5780049385015"
$ws.Range("C16").Value = "400.0000001247187"
$ws.Range("I16").Value = "125.61235673704185"
$ws.Range("C17").Value = "400.0000001247187"
$ws.Range("I17").Value = "-84.50830206920877"
$ws.Range("C18").Value = "400.0000001250004"
$ws.Range("I18").Value = "253.68208029388555"
$ws.Range("C19").Value = "400.00000012471844"
$ws.Range("I19").Value = "-580.0586547098526"
$ws.Range("I20").Value = "126.60289474404493"
$ws.Range("C21").Value = "170.39999387651875"
$ws.Range("I21").Value = "-385.26448161475287"
$ws.Range("C22").Value = "400.00000011271607"
$ws.Range("I22").Value = "-732.3309334884935"
$ws.Range("C23").Value = "400.000000119644"
$ws.Range("I23").Value = "-806.7850021760776"
$ws.Range("I24").Value = "-372.3030611317588"
$ws.Range("C25").Value = "420.00000009118696"
$ws.Range("C26").Value = "491.8359989661517"
$ws.Range("I26").Value = "37.33103371373209"
$ws.Range("C27").Value = "400.000000119644"
$ws.Range("I27").Value = "-119.22099286304184"
$ws.Range("I28").Value = "78.88041297218238"
$ws.Range("I29").Value = "-70.46312705073609"
$ws.Range("I30").Value = "-70.46312705073609"
$ws.Range("I31").Value = "-157.49641525530498"
$ws.Range("I32").Value = "-157.49641525530498"
$ws.Range("I33").Value = "202.8786633189661"
$ws.Range("I34").Value = "202.87866331896646"
$ws.Range("I35").Value = "-56.346568442022885"
$ws.Range("I36").Value = "-56.34656844202289"

$ws = $wb.Worksheets.Item(8)
$ws.Range("C2").Value = "83.1102981641134"
$ws.Range("I2").Value = "-43.76604388854654"
$ws.Range("C3").Value = "400.0000000176049"
$ws.Range("I3").Value = "-462.8605050175178"
$ws.Range("C4").Value = "400.0000000356312"
$ws.Range("I4").Value = "-389.3696336080531"
$ws.Range("C5").Value = "420.0000000374189"
$ws.Range("I5").Value = "-404.198895308477"
$ws.Range("C6").Value = "457.0661199811953"
$ws.Range("I6").Value = "78.23974826787152"
$ws.Range("C7").Value = "304.7895961264533"
$ws.Range("I7").Value = "37.449797108514815"
$ws.Range("C8").Value = "299.3835951197538"
$ws.Range("I8").Value = "-119.15139987898395"
$ws.Range("C9").Value = "400.0000000410847"
$ws.Range("I9").Value = "-98.59411179717806"
$ws.Range("C10").Value = "400.00000004108404"
$ws.Range("I10").Value = "700.066452836654"
$ws.Range("J10").Value = "796.7433714822"
$ws.Range("C11").Value = "400.00000004041436"
$ws.Range("I11").Value = "1263.936199157161"
$ws.Range("C12").Value = "400.0000000412478"
$ws.Range("I12").Value = "1125.7950304517337"
$ws.Range("C13").Value = "400.0000000410847"
$ws.Range("I13").Value = "714.0676772664089"
$ws.Range("C14").Value = "400.0000000412478"
$ws.Range("C15").Value = "492.0673117782611"
$ws.Range("I15").Value = "958.9250499505918"
$ws.Range("C16").Value = "400.0000001247187"
$ws.Range("I16").Value = "-83.94099201703399"
$ws.Range("C17").Value = "400.0000001247187"
$ws.Range("I17").Value = "-247.86327821077228"
$ws.Range("C18").Value = "400.0000001250004"
$ws.Range("I18").Value = "253.81990997037437"
$ws.Range("C19").Value = "400.00000012471844"
$ws.Range("I19").Value = "-579.717058598511"
$ws.Range("C20").Value = "500.0"
$ws.Range("I20").Value = "-82.87517660503312"
$ws.Range("C21").Value = "170.39999387651875"
$ws.Range("I21").Value = "-549.5287560536098"
$ws.Range("I22").Value = "-732.3876796489341"
$ws.Range("C23").Value = "400.000000119644"
$ws.Range("I23").Value = "-806.9480828530834"
$ws.Range("C24").Value = "399.99999994878533"
$ws.Range("I24").Value = "-583.8855436010215"
$ws.Range("C25").Value = "420.00000009118696"
$ws.Range("I25").Value = "78.2397482678715"
$ws.Range("C26").Value = "492.0750428452021"
$ws.Range("C27").Value = "400.000000119644"
$ws.Range("I27").Value = "-119.15139987898398"
$ws.Range("I29").Value = "-75.85131140530805"
$ws.Range("I30").Value = "-75.85131140530805"
$ws.Range("I31").Value = "-153.40200364548582"
$ws.Range("I32").Value = "-153.4020036454986"
$ws.Range("I33").Value = "203.16659535495685"
$ws.Range("I34").Value = "203.1665953549572"
$ws.Range("I35").Value = "-62.017703835741536"
$ws.Range("I36").Value = "-62.017703835741536"

$ws = $wb.Worksheets.Item(9)
$ws.Range("C2").Value = "83.1102981641134"
$ws.Range("C3").Value = "400.0000000176049"
$ws.Range("I3").Value = "-463.0886366549243"
$ws.Range("C4").Value = "400.0000000356312"
$ws.Range("I4").Value = "-389.67697253383494"
$ws.Range("C5").Value = "420.0000000374189"
$ws.Range("I5").Value = "-261.5957845235898"
$ws.Range("C6").Value = "455.40363784850354"
$ws.Range("I6").Value = "100.0473134066576"
$ws.Range("C7").Value = "333.74967258498094"
$ws.Range("C8").Value = "299.3835951197538"
$ws.Range("I8").Value = "-119.35724061259066"
$ws.Range("C9").Value = "400.0000000410847"
$ws.Range("C10").Value = "400.00000004108404"
$ws.Range("I10").Value = "654.7535952677578"
$ws.Range("J10").Value = "796.7433714822"
$ws.Range("C11").Value = "400.00000004041436"
$ws.Range("I11").Value = "1262.1643523836742"
$ws.Range("C12").Value = "400.0000000412478"
$ws.Range("I12").Value = "1124.5523333067279"
$ws.Range("C13").Value = "400.0000000410847"
$ws.Range("I13").Value = "939.9893068077228"
$ws.Range("C14").Value = "400.0000000412478"
$ws.Range("I14").Value = "-44.916430386998826"
$ws.Range("C15").Value = "490.40453435910837"
$ws.Range("I15").Value = "957.5504434202403"
$ws.Range("C16").Value = "400.0000001247187"
$ws.Range("I16").Value = "151.03782733273124"
$ws.Range("C17").Value = "400.0000001247187"
$ws.Range("I17").Value = "-251.1920171128559"
$ws.Range("C18").Value = "400.0000001250004"
$ws.Range("I18").Value = "253.2785196352818"
$ws.Range("C19").Value = "400.00000012471844"
$ws.Range("I19").Value = "-579.8742404133734"
$ws.Range("C20").Value = "500.0"
$ws.Range("I20").Value = "-230.71222905433567"
$ws.Range("C21").Value = "170.39999387651875"
$ws.Range("I21").Value = "-552.6173380860047"
$ws.Range("C22").Value = "400.00000011271607"
$ws.Range("I22").Value = "-732.1626886104333"
$ws.Range("I23").Value = "-806.6415234857512"
$ws.Range("C24").Value = "399.99999994878533"
$ws.Range("I24").Value = "-726.055595930467"
$ws.Range("C25").Value = "420.00000009118696"
$ws.Range("I25").Value = "100.04731340665762"
$ws.Range("C26").Value = "429.2542886853422"
$ws.Range("I26").Value = "36.98169615132633"
$ws.Range("C27").Value = "400.000000119644"
$ws.Range("I27").Value = "-119.35724061259066"
$ws.Range("I29").Value = "-139.3797113007853"
$ws.Range("I31").Value = "-173.52200328373493"
$ws.Range("I32").Value = "-173.52200328373493"
$ws.Range("I33").Value = "202.02650413115003"
$ws.Range("I34").Value = "202.02650413115043"
$ws.Range("I36").Value = "-111.52426486321238"
